$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.359.20'
$ws.Range("E2").Value = '  +1.56%  '

$ws.Range("D3").Value = '2.638.46'
$ws.Range("E3").Value = '  +1.41%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '599.59'
$ws.Range("E5").Value = '  +1.13%  '

$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '154.72'
$ws.Range("E6").Value = '  +3.10%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '0.545'
$ws.Range("E8").Value = '  +0.54%  '

$ws.Range("D9").Value = '2.637.63'
$ws.Range("E9").Value = '  +1.37%  '

$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '0.138'
$ws.Range("E10").Value = '  +2.46%  '

$ws.Range("E11").Value = '  -0.39%  '

$ws.Range("E12").Value = '  +1.18%  '

$ws.Range("E13").Value = '  +2.17%  '

$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '28.10'
$ws.Range("E14").Value = '  +3.35%  '

$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '0.0000187'
$ws.Range("E15").Value = '  +0.92%  '

$ws.Range("D16").Value = '3.116.75'
$ws.Range("E16").Value = '  +1.29%  '

$ws.Range("D17").Value = '68.254.24'
$ws.Range("E17").Value = '  +1.52%  '

$ws.Range("D18").Value = '2.645.67'
$ws.Range("E18").Value = '  +1.59%  '

$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '11.51'
$ws.Range("E19").Value = '  +4.65%  '

$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '365.89'
$ws.Range("E20").Value = '  -1.26%  '

$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '7.48'
$ws.Range("E21").Value = '  +1.84%  '

$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '4.34'
$ws.Range("E22").Value = '  +3.70%  '

$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '4.88'
$ws.Range("E23").Value = '  +2.78%  '

$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '2.11'
$ws.Range("E24").Value = '  +4.48%  '

$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '73.29'
$ws.Range("E25").Value = '  +10.02%  '

$ws.Range("E26").Value = '  -0.01%  '

$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '9.97'
$ws.Range("E27").Value = '  +1.16%  '

$ws.Range("B28").Value = 'WrappedeETH'
$ws.Range("C28").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D28").Value = '2.774.92'
$ws.Range("E28").Value = '  +1.37%  '

$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '0.0000105'
$ws.Range("E29").Value = '  +5.56%  '

$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.27%  '

$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '573.32'
$ws.Range("E31").Value = '  -1.37%  '

$ws.Range("E32").Value = '  +4.84%  '

$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '8.03'
$ws.Range("E33").Value = '  +4.75%  '

$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '1.87'
$ws.Range("E34").Value = '  +3.01%  '

$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '0.130'

$ws.Range("E36").Value = '  +0.04%  '

$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '1.57'
$ws.Range("E37").Value = '  +4.48%  '

$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '160.45'
$ws.Range("E38").Value = '  +1.97%  '

$ws.Range("E39").Value = '  +1.71%  '

$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '1.91'
$ws.Range("E40").Value = '  +2.84%  '

$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '0.371'
$ws.Range("E41").Value = '  +1.54%  '

$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '5.42'
$ws.Range("E42").Value = '  +4.32%  '

$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '2.68'
$ws.Range("E43").Value = '  +4.79%  '

$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '17.74'
$ws.Range("E44").Value = '  +3.76%  '

$ws.Range("D45").Value = '0.0₆0320'
$ws.Range("E45").Value = '  +13.51%  '

$ws.Range("E46").Value = '  +0.05%  '

$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '40.47'
$ws.Range("E47").Value = '  +0.82%  '

$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '157.83'
$ws.Range("E48").Value = '  +3.58%  '

$ws.Range("E49").Value = '  +3.77%  '

$ws.Range("E50").Value = '  +2.95%  '

$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '22.02'
$ws.Range("E51").Value = '  +4.18%  '
